# This script updates the "想去人数" (interest count) column F values
# on the "展览" and "全部类型" worksheets to reflect newly generated
# output data (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 2029
$ws1.Range("F6").Value  = 575
$ws1.Range("F8").Value  = 2058
$ws1.Range("F9").Value  = 10473
$ws1.Range("F11").Value = 152
$ws1.Range("F12").Value = 274
$ws1.Range("F15").Value = 7360
$ws1.Range("F18").Value = 182
$ws1.Range("F20").Value = 3289

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 2029
$ws4.Range("F6").Value  = 575
$ws4.Range("F9").Value  = 2058
$ws4.Range("F12").Value = 10473
$ws4.Range("F14").Value = 152
$ws4.Range("F15").Value = 274
$ws4.Range("F18").Value = 7360
$ws4.Range("F21").Value = 182
$ws4.Range("F23").Value = 3289
